$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6; everything from the old row 6 downward
# (through the old row 63) shifts down to rows 7-64, carrying its
# formatting (incl. the date number format on column D) along with it.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly record.
$ws.Cells.Item(6, 1).Value = 8
$ws.Cells.Item(6, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(6, 3).Value = "Coquimbo"
$ws.Cells.Item(6, 4).Value = 44530
$ws.Cells.Item(6, 5).Value = 4
$ws.Cells.Item(6, 6).Value = 100112052
$ws.Cells.Item(6, 7).Value = "Albahaca"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 800
$ws.Cells.Item(6, 11).Value = 3000
$ws.Cells.Item(6, 12).Value = 4000
$ws.Cells.Item(6, 13).Value = 3500
$ws.Cells.Item(6, 14).Value = "`$/paquete"
$ws.Cells.Item(6, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(6, 16).Value = 3500
$ws.Cells.Item(6, 17).Value = 1
$ws.Cells.Item(6, 18).Value = "Hortaliza"
